# Weekly Fruit/Vegetable data update for "Ajo" (Garlic) sheet.
# Two new daily records are inserted at the top of the data block (rows 98-99),
# pushing all the existing records (old rows 98-124) down by two rows
# (new rows 100-126). The workbook dimension grows from A1:R124 to A1:R126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 98, shifting
# everything from row 98 down onward two rows further down.
$ws.Range("A98:A99").EntireRow.Insert()

# --- New row 98 ---
$ws.Range("A98").Value = 9
$ws.Range("B98").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C98").Value = "Metropolitana"
$ws.Range("D98").Value = 44508
$ws.Range("E98").Value = 13
$ws.Range("F98").Value = 100112003
$ws.Range("G98").Value = "Ajo"
$ws.Range("H98").Value = "Chino"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 430
$ws.Range("K98").Value = 17000
$ws.Range("L98").Value = 18000
$ws.Range("M98").Value = 17500
$ws.Range("N98").Value = "$/caja 10 kilos"
$ws.Range("O98").Value = "China"
$ws.Range("P98").Value = 1750
$ws.Range("Q98").Value = 10
$ws.Range("R98").Value = "Hortaliza"

# --- New row 99 ---
$ws.Range("A99").Value = 9
$ws.Range("B99").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C99").Value = "Metropolitana"
$ws.Range("D99").Value = 44508
$ws.Range("E99").Value = 13
$ws.Range("F99").Value = 100112003
$ws.Range("G99").Value = "Ajo"
$ws.Range("H99").Value = "Chino"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 340
$ws.Range("K99").Value = 16000
$ws.Range("L99").Value = 17000
$ws.Range("M99").Value = 16500
$ws.Range("N99").Value = "$/malla 10 kilos"
$ws.Range("O99").Value = "China"
$ws.Range("P99").Value = 1650
$ws.Range("Q99").Value = 10
$ws.Range("R99").Value = "Hortaliza"
